$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.856.84"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "2.120.93"
$ws.Range("E3").Value = "  +10.74%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "256.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.672"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.17%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.40"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.62"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.375"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0743"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.40%  "
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").Value = "2.431.72"
$ws.Range("E13").Value = "  +10.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.853"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.88%  "
$ws.Range("D16").Value = "2.123.09"
$ws.Range("E16").Value = "  +10.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").Value = "36.920.77"
$ws.Range("E18").Value = "  +1.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "0.0₃0846"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "242.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.70%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  -8.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "172.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +13.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.10%  "
$ws.Range("E29").Value = "  -6.53%  "
$ws.Range("E30").Value = "  -3.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +51.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0951"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +12.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0603"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +19.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.22"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.20%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.922"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.12%  "
$ws.Range("E41").Value = "  +8.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0225"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +16.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.23%  "
$ws.Range("D46").Value = "1.365.48"
$ws.Range("E46").Value = "  +1.85%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +12.05%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0841"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.15%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.313.56"
$ws.Range("E49").Value = "  +10.56%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("E51").Value = "  +1.65%  "
